$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -------------------------------------------------------
# Drop the two extra "object" rows (3 and 4) - only two rows remain.
$ws.Rows("3:4").Delete()

# Fix the capitalisation / rename the entries for the remaining rows.
$ws.Range("A1").Value = "Teste"
$ws.Range("A2").Value = "Tv da sala"
$ws.Range("B2").Value = "Televisor"

# Row 2 no longer carries a temperature (C2) and its "on" flag moves from
# D2 (boolean) to a plain numeric 0, with a new boolean flag in E2.
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = $false

# --- Column layout / formatting ---------------------------------------
# Columns A and B: plain (general) text columns, resized to fit.
$ws.Columns("A:B").ColumnWidth = 12.43357142857143
$ws.Columns("A:B").HorizontalAlignment = 1

# Column C: numeric column, right aligned with a thousands separator.
$ws.Columns("C").ColumnWidth = 12.43357142857143
$ws.Columns("C").NumberFormat = "#,##0"
$ws.Columns("C").HorizontalAlignment = -4152

# Column D: boolean "button" column, centered.
$ws.Columns("D").ColumnWidth = 12.43357142857143
$ws.Columns("D").HorizontalAlignment = -4108
